$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (A11=9): description stays "blend 4,5 with sigmoid_drop using lasso_drop"
# (the old "...using lasso" entry, previously in row 12, is being removed entirely)
# Row 12 (A12=10): is repurposed for the new "lasso_minimal_14vars" trial; its CV
# value is cleared since it hasn't been scored yet.
# New column M holds expectation notes for rows 10-12.
# Order of assignment matters for shared-string table ordering on save.
$ws.Range("M11").Value = "Expect this to beat 8"
$ws.Range("M10").Value = "Expect this to beat 7"
$ws.Range("B12").Value = "lasso_minimal_14vars"
$ws.Range("D12").ClearContents()
$ws.Range("M12").Value = "Expect this to beat 9"

# Update selection to reflect where the author ended up working
[void]$ws.Range("M12").Select()
